$d = $word.ActiveDocument

# Locate the "Dewalt 20 Volt Power Tools" table by finding the row whose
# first cell contains "DCB119" (the current last row of that table).
$target = $null
foreach ($t in $d.Tables) {
    $lastRow = $t.Rows.Item($t.Rows.Count)
    $cellText = $lastRow.Cells.Item(1).Range.Text.TrimEnd([char]13, [char]7)
    if ($cellText -eq "DCB119") {
        $target = $t
        break
    }
}

$newRows = @(
    @("DCPR320B", "20V MAX* 1-1/2 in Cordless Pruner"),
    @("DCM848P2", "20V MAX* XR® 5 in. Cordless Variable-Speed Random Orbit Polisher Kit"),
    @("DCM849P2", "20V MAX* XR® 7 in. Cordless Variable-Speed Rotary Polisher Kit")
)

foreach ($pair in $newRows) {
    $row = $target.Rows.Add()
    $row.Cells.Item(1).Range.Text = $pair[0]
    $row.Cells.Item(2).Range.Text = $pair[1]
}
